# Refresh the cryptocurrency price/volume snapshot (cryptos.xlsx).
# Updates the "Price" (D) and "Volume(1h)" (E) columns for rows 2-51
# with the latest scraped values, as produced by the scheduled
# GitHub Actions job that regenerates this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking values (e.g. "1.000",
# "315.67") as plain text in the source data. Force those cells to
# text format before writing so Excel keeps the literal string
# (incl. trailing zeros / thousands-dot formatting) instead of
# silently re-interpreting it as a number.
$priceCells = @(
    "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9",
    "D10", "D13", "D15", "D16", "D17", "D19", "D21", "D23",
    "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D32",
    "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40",
    "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48",
    "D49", "D50"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.932.60"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "1.777.56"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "315.67"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.5377"
$ws.Range("E7").Value = "  -2.48%  "
$ws.Range("D8").Value = "0.3764"
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").Value = "0.07436"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").Value = "41.62"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("E11").Value = "  -2.61%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "20.42"
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").Value = "7.193"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "1.774.74"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").Value = "87.98"
$ws.Range("E17").Value = "  -4.76%  "
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").Value = "0.06425"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "17.19"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("D23").Value = "27.960.67"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").Value = "11.12"
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").Value = "2.080"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").Value = "155.85"
$ws.Range("D27").Value = "20.21"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "1.972.30"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("D29").Value = "2.271"
$ws.Range("E29").Value = "  -5.45%  "
$ws.Range("D30").Value = "119.69"
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").Value = "0.1052"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").Value = "3.641"
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("D34").Value = "5.504"
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("D35").Value = "0.2250"
$ws.Range("D36").Value = "0.06369"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").Value = "0.02261"
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("D38").Value = "4.957"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("D39").Value = "8.384"
$ws.Range("E39").Value = "  -5.90%  "
$ws.Range("D40").Value = "0.6117"
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("D41").Value = "11.00"
$ws.Range("E41").Value = "  -5.39%  "
$ws.Range("D42").Value = "1.177"
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").Value = "1.427"
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "13.19"
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").Value = "3.653"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "0.5729"
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("D48").Value = "126.21"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "1.184"
$ws.Range("E49").Value = "  +3.11%  "
$ws.Range("D50").Value = "1.919"
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("E51").Value = "  -1.86%  "
